$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Cells changing data type (numeric <-> shared-string text "0" / "***.*") ---
# Copying from an existing cell that already holds the exact target value/style
# preserves both the correct value and the correct style/number-format.
$ws.Range("M14").Copy($ws.Range("N14"))   # -100 (numeric, style 16)
$ws.Range("C14").Copy($ws.Range("D15"))   # "0"   (text, style 14)
$ws.Range("E14").Copy($ws.Range("E15"))   # "***.*" (text, style 14)
$ws.Range("C14").Copy($ws.Range("C20"))   # "0"
$ws.Range("C14").Copy($ws.Range("D22"))   # "0"
$ws.Range("E14").Copy($ws.Range("E22"))   # "***.*"
$ws.Range("C14").Copy($ws.Range("C27"))   # "0"
$ws.Range("C14").Copy($ws.Range("D30"))   # "0"
$ws.Range("E14").Copy($ws.Range("E30"))   # "***.*"

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = -25
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -76.923076923076

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -36.363636363636
$ws.Range("I16").Value = 16
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = 6.666666666666
$ws.Range("L16").Value = -27.272727272727
$ws.Range("M16").Value = -57.894736842105
$ws.Range("N16").Value = -93.650793650793

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -19.354838709677
$ws.Range("I17").Value = 43
$ws.Range("J17").Value = 48
$ws.Range("K17").Value = -10.416666666666
$ws.Range("L17").Value = 2.380952380952
$ws.Range("M17").Value = 2.380952380952
$ws.Range("N17").Value = -65.322580645161

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 21
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 5
$ws.Range("L18").Value = -16
$ws.Range("M18").Value = -53.333333333333
$ws.Range("N18").Value = -86.708860759493

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 26.923076923076
$ws.Range("I19").Value = 48
$ws.Range("J19").Value = 42
$ws.Range("K19").Value = 14.285714285714
$ws.Range("L19").Value = 37.142857142857
$ws.Range("N19").Value = -49.473684210526

# --- Row 20 (G.L.A.) ---
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 75
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = -35.294117647058
$ws.Range("N20").Value = -86.25

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 3.409090909090
$ws.Range("I21").Value = 142
$ws.Range("J21").Value = 141
$ws.Range("K21").Value = 0.709219858156
$ws.Range("L21").Value = 3.649635036496
$ws.Range("M21").Value = -21.978021978022
$ws.Range("N21").Value = -80.413793103448

# --- Row 23 (Housing) ---
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 13
$ws.Range("H23").Value = -18.75
$ws.Range("I23").Value = 26
$ws.Range("J23").Value = 29
$ws.Range("K23").Value = -10.344827586206
$ws.Range("L23").Value = 8.333333333333
$ws.Range("M23").Value = 85.714285714285

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -6.666666666666
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 110
$ws.Range("H24").Value = 3.636363636363
$ws.Range("I24").Value = 178
$ws.Range("J24").Value = 166
$ws.Range("K24").Value = 7.228915662650
$ws.Range("L24").Value = 37.984496124031
$ws.Range("M24").Value = 69.523809523809

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -36.363636363636
$ws.Range("F25").Value = 43
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = -20.370370370370
$ws.Range("I25").Value = 79
$ws.Range("J25").Value = 76
$ws.Range("K25").Value = 3.947368421052
$ws.Range("L25").Value = 68.085106382978
$ws.Range("M25").Value = -8.139534883720

# --- Row 26 (UCR Rape*) ---
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 6
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 20

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 300
